$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1302.3914
$ws.Range("I40").Value = 1150.3846
$ws.Range("J40").Value = 1500
$ws.Range("K40").Value = 1150.3846
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = -975.3846000000001
$ws.Range("N40").Value = -1850
$ws.Range("H64").Value = 3351.7856
$ws.Range("I64").Value = 2927.1428
$ws.Range("K64").Value = 2927.1428
$ws.Range("M64").Value = -2679.1428
$ws.Range("H67").Value = 3351.7856
$ws.Range("I67").Value = 2927.1428
$ws.Range("K67").Value = 2927.1428
$ws.Range("M67").Value = -2069.1428
$ws.Range("H138").Value = 1981.74
$ws.Range("I138").Value = 901.2857
$ws.Range("J138").Value = 2157.628
$ws.Range("K138").Value = 2703.8571
$ws.Range("L138").Value = 6472.884
$ws.Range("M138").Value = 2436.1429
$ws.Range("N138").Value = -16752.884

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 704.5568
$ws.Range("I31").Value = 635.0345
$ws.Range("J31").Value = 838.9666999999999
$ws.Range("K31").Value = 635.0345
$ws.Range("L31").Value = 838.9666999999999
$ws.Range("M31").Value = -340.0345
$ws.Range("N31").Value = -1428.9667
$ws.Range("H34").Value = 704.5568
$ws.Range("I34").Value = 635.0345
$ws.Range("J34").Value = 838.9666999999999
$ws.Range("K34").Value = 635.0345
$ws.Range("L34").Value = 838.9666999999999
$ws.Range("M34").Value = -433.0345
$ws.Range("N34").Value = -1242.9667
$ws.Range("H62").Value = 40002600
$ws.Range("I62").Value = 2666.3333
$ws.Range("J62").Value = 100002504
$ws.Range("K62").Value = 2666.3333
$ws.Range("L62").Value = 100002504
$ws.Range("M62").Value = -2042.3333
$ws.Range("N62").Value = -100003752
$ws.Range("H65").Value = 40002600
$ws.Range("I65").Value = 2666.3333
$ws.Range("J65").Value = 100002504
$ws.Range("K65").Value = 13331.6665
$ws.Range("L65").Value = 500012520
$ws.Range("M65").Value = -10211.6665
$ws.Range("N65").Value = -500018760

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 4048.4614
$ws.Range("J64").Value = 4762.8
$ws.Range("L64").Value = 14288.4
$ws.Range("N64").Value = -14828.4
$ws.Range("H67").Value = 4048.4614
$ws.Range("J67").Value = 4762.8
$ws.Range("L67").Value = 14288.4
$ws.Range("N67").Value = -16160.4
$ws.Range("H76").Value = 6163.433
$ws.Range("I76").Value = 3966.6667
$ws.Range("J76").Value = 6407.5186
$ws.Range("K76").Value = 11900.0001
$ws.Range("L76").Value = 19222.5558
$ws.Range("M76").Value = -11517.0001
$ws.Range("N76").Value = -19988.5558
$ws.Range("H79").Value = 6163.433
$ws.Range("I79").Value = 3966.6667
$ws.Range("J79").Value = 6407.5186
$ws.Range("K79").Value = 11900.0001
$ws.Range("L79").Value = 19222.5558
$ws.Range("M79").Value = -10574.0001
$ws.Range("N79").Value = -21874.5558
$ws.Range("H88").Value = 3833.3333
$ws.Range("J88").Value = 8083.3335
$ws.Range("L88").Value = 24250.0005
$ws.Range("N88").Value = -25106.0005
$ws.Range("H91").Value = 3833.3333
$ws.Range("J91").Value = 8083.3335
$ws.Range("L91").Value = 24250.0005
$ws.Range("N91").Value = -27214.0005
$ws.Range("H94").Value = 4685
$ws.Range("I94").Value = 1975
$ws.Range("J94").Value = 5362.5
$ws.Range("K94").Value = 5925
$ws.Range("L94").Value = 16087.5
$ws.Range("M94").Value = -5249
$ws.Range("N94").Value = -17439.5
$ws.Range("H97").Value = 237.45454
$ws.Range("I97").Value = 115.42857
$ws.Range("J97").Value = 451
$ws.Range("K97").Value = 346.28571
$ws.Range("L97").Value = 1353
$ws.Range("M97").Value = 149.71429
$ws.Range("N97").Value = -2345
$ws.Range("H100").Value = 3207.5557
$ws.Range("J100").Value = 3207.5557
$ws.Range("L100").Value = 9622.667099999999
$ws.Range("N100").Value = -11244.6671
$ws.Range("H103").Value = 3054.6365
$ws.Range("I103").Value = 918.25
$ws.Range("J103").Value = 4275.4287
$ws.Range("K103").Value = 2754.75
$ws.Range("L103").Value = 12826.2861
$ws.Range("M103").Value = -1875.75
$ws.Range("N103").Value = -14584.2861
$ws.Range("H106").Value = 4119.846
$ws.Range("J106").Value = 4119.846
$ws.Range("L106").Value = 12359.538
$ws.Range("N106").Value = -14251.538
$ws.Range("H109").Value = 79210.30499999999
$ws.Range("I109").Value = 125966.75
$ws.Range("J109").Value = 4400
$ws.Range("K109").Value = 377900.25
$ws.Range("L109").Value = 13200
$ws.Range("M109").Value = -376860.25
$ws.Range("N109").Value = -15280
$ws.Range("H112").Value = 90923460
$ws.Range("J112").Value = 100015400
$ws.Range("L112").Value = 300046200
$ws.Range("N112").Value = -300048416
$ws.Range("H115").Value = 4675.1333
$ws.Range("I115").Value = 2671.1667
$ws.Range("J115").Value = 6011.1113
$ws.Range("K115").Value = 8013.500100000001
$ws.Range("L115").Value = 18033.3339
$ws.Range("M115").Value = -6838.500100000001
$ws.Range("N115").Value = -20383.3339
$ws.Range("H118").Value = 510
$ws.Range("I118").Value = 510
$ws.Range("K118").Value = 1530
$ws.Range("M118").Value = -287
$ws.Range("H121").Value = 608.25
$ws.Range("I121").Value = 247.42857
$ws.Range("J121").Value = 888.8889
$ws.Range("K121").Value = 742.28571
$ws.Range("L121").Value = 2666.6667
$ws.Range("M121").Value = 567.71429
$ws.Range("N121").Value = -5286.6667

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = $null
$ws.Range("N47").Value = $null
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = $null
$ws.Range("N52").Value = $null
$ws.Range("H136").Value = 4833
$ws.Range("I136").Value = 4833
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 14499
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -11949
$ws.Range("N136").Value = $null

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2115.2856
$ws.Range("I136").Value = 1949.5
$ws.Range("J136").Value = 2181.6
$ws.Range("K136").Value = 5848.5
$ws.Range("L136").Value = 6544.799999999999
$ws.Range("M136").Value = -3298.5
$ws.Range("N136").Value = -11644.8
